$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "ID" values up so the new hire lands at the bottom of the list ---
# Before: C34=1C106BD2, C35=C4090B29, C36=A46E6932 (no C37)
# After:  C17=1C106BD2, C36=C4090B29, C37=A46E6932 (C34/C35 cleared)
$ws.Range("C17").Value = "1C106BD2"
$ws.Range("C34").ClearContents()
$ws.Range("C35").ClearContents()
$ws.Range("C36").Value = "C4090B29"

# --- Add the new worker row ---
$ws.Range("A37").Value = "ZyAsia Holmes"
$ws.Range("B37").Value = 1185
$ws.Range("C37").Value = "A46E6932"

# Highlight the newly-added row's name cell with a themed border
# (Green, Accent 6, Lighter 40%) on the right/top/bottom edges only.
$newRow = $ws.Range("A37")
$newRow.Borders.Color = 9359785
$newRow.Borders.Item(7).LineStyle = 0

# --- Update the active selection to where the user was last working ---
$ws.Range("D20").Select() | Out-Null
